# Update gh-pages to output generated at 456a3b4
#
# Refreshes the "want to go" counters (column F) for several existing
# events and inserts a new "南宁·国乙only" row (2024-08-10) into the two
# worksheets that list exhibition-type events: 展览 (sheet 1) and
# 全部类型 (sheet 4). The 演出 / 本地生活 sheets are untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 597
$ws1.Range("F4").Value = 453
$ws1.Range("F5").Value = 485
$ws1.Range("F6").Value = 278
$ws1.Range("F7").Value = 2533
$ws1.Range("F9").Value = 6807
$ws1.Range("F10").Value = 180
$ws1.Range("F11").Value = 432

# Insert the new event as row 12, pushing "万圣漫控嘉年华10" down to row 13.
$ws1.Rows.Item(12).Insert()

# Copy A11's style to A12 (keeps the bordered/centered look of column A),
# then overwrite the value with the correct sequence number.
$ws1.Range("A11").Copy($ws1.Range("A12"))
$ws1.Range("A12").Value = 11

# Force text storage so the date-looking string isn't auto-converted into
# an Excel date serial.
$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = "2024-08-10"
$ws1.Range("C12").Value = "南宁·国乙only"
$ws1.Range("D12").Value = "新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店"
$ws1.Range("E12").Value = "2024.08.10 10:00-08.10 17:00"
$ws1.Range("F12").Value = 1
$ws1.Range("G12").Value = 40
$ws1.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=88227"
$ws1.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202406/3cFX9LLQ1719482186347.jpeg"

# Row 13 is the old row 12 shifted down by the insert; refresh its index
# and "want to go" counter.
$ws1.Range("A13").Value = 12
$ws1.Range("F13").Value = 32

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 597
$ws4.Range("F4").Value = 453
$ws4.Range("F5").Value = 485
$ws4.Range("F6").Value = 278
$ws4.Range("F9").Value = 2533
$ws4.Range("F11").Value = 6807
$ws4.Range("F12").Value = 180
$ws4.Range("F13").Value = 432

# Insert the new event as row 14, pushing the remaining rows down by one.
$ws4.Rows.Item(14).Insert()

$ws4.Range("A13").Copy($ws4.Range("A14"))
$ws4.Range("A14").Value = 13

$ws4.Range("B14").NumberFormat = "@"
$ws4.Range("B14").Value = "2024-08-10"
$ws4.Range("C14").Value = "南宁·国乙only"
$ws4.Range("D14").Value = "新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店"
$ws4.Range("E14").Value = "2024.08.10 10:00-08.10 17:00"
$ws4.Range("F14").Value = 1
$ws4.Range("G14").Value = 40
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=88227"
$ws4.Range("I14").Value = "//i0.hdslb.com/bfs/openplatform/202406/3cFX9LLQ1719482186347.jpeg"

# Rows 15-17 are the old rows 14-16 shifted down by the insert; refresh
# their sequence numbers and the trailing "万圣漫控嘉年华10" counter.
$ws4.Range("A15").Value = 14
$ws4.Range("A16").Value = 15
$ws4.Range("A17").Value = 16
$ws4.Range("F17").Value = 32

Write-Output "Applied gh-pages refresh (456a3b4)"
